$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.017573522171699
$ws.Range("D2").Value = 1.032216671736316
$ws.Range("E2").Value = 1.01895412357485
$ws.Range("F2").Value = 1.015937847442841
$ws.Range("I2").Value = 1.031487194878901
$ws.Range("J2").Value = 1.022786867894597
$ws.Range("K2").Value = 1.035022541285615
$ws.Range("L2").Value = 1.021798806911659
$ws.Range("M2").Value = 1.018791512133567
$ws.Range("N2").Value = 1.011687332015984
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.018431927234231
$ws.Range("D3").Value = 1.032702661840093
$ws.Range("E3").Value = 1.019679667576192
$ws.Range("F3").Value = 1.017436600937688
$ws.Range("I3").Value = 1.031620742544686
$ws.Range("J3").Value = 1.023281697063658
$ws.Range("K3").Value = 1.035318121683775
$ws.Range("L3").Value = 1.02233049326863
$ws.Range("M3").Value = 1.020093617021657
$ws.Range("N3").Value = 1.011852679580706
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.01898757080724
$ws.Range("D4").Value = 1.033017203349916
$ws.Range("E4").Value = 1.020149708460257
$ws.Range("F4").Value = 1.018406564735827
$ws.Range("I4").Value = 1.03170600099857
$ws.Range("J4").Value = 1.02360146412114
$ws.Range("K4").Value = 1.035508734712094
$ws.Range("L4").Value = 1.022674429038855
$ws.Range("M4").Value = 1.020935851064403
$ws.Range("N4").Value = 1.011959482729608
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.019221210272095
$ws.Range("D5").Value = 1.033149452501587
$ws.Range("E5").Value = 1.020347448092494
$ws.Range("F5").Value = 1.018814382505304
$ws.Range("I5").Value = 1.031741566469893
$ws.Range("J5").Value = 1.023735793179065
$ws.Range("K5").Value = 1.035588712409683
$ws.Range("L5").Value = 1.022818994956824
$ws.Range("M5").Value = 1.021289854454292
$ws.Range("N5").Value = 1.012004337647208
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.019260442101867
$ws.Range("D6").Value = 1.033171658599658
$ws.Range("E6").Value = 1.020380657304102
$ws.Range("F6").Value = 1.018882859665636
$ws.Range("I6").Value = 1.031747521800986
$ws.Range("J6").Value = 1.023758341685476
$ws.Range("K6").Value = 1.035602131836898
$ws.Range("L6").Value = 1.022843266747728
$ws.Range("M6").Value = 1.021349289074778
$ws.Range("N6").Value = 1.012011866334033
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.018990692525544
$ws.Range("D7").Value = 1.033018970409196
$ws.Range("E7").Value = 1.020152350139783
$ws.Range("F7").Value = 1.018412013834721
$ws.Range("I7").Value = 1.031706477315876
$ws.Range("J7").Value = 1.023603259429936
$ws.Range("K7").Value = 1.035509803992401
$ws.Range("L7").Value = 1.022676360833381
$ws.Range("M7").Value = 1.020940581552696
$ws.Range("N7").Value = 1.011960082260745
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.017863582853292
$ws.Range("D8").Value = 1.032380898161505
$ws.Range("E8").Value = 1.019199206949922
$ws.Range("F8").Value = 1.016444325562885
$ws.Range("I8").Value = 1.031532566970073
$ws.Range("J8").Value = 1.022954184591376
$ws.Range("K8").Value = 1.035122567444039
$ws.Range("L8").Value = 1.021978513367785
$ws.Range("M8").Value = 1.019231631934139
$ws.Range("N8").Value = 1.011743250778884
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.01587901125328
$ws.Range("D9").Value = 1.031257169330683
$ws.Range("E9").Value = 1.017524017094466
$ws.Range("F9").Value = 1.012978118701747
$ws.Range("I9").Value = 1.031217283563179
$ws.Range("J9").Value = 1.021807226035386
$ws.Range("K9").Value = 1.034435296853642
$ws.Range("L9").Value = 1.020748067794838
$ws.Range("M9").Value = 1.016217672278221
$ws.Range("N9").Value = 1.011359735222549
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.014557026781007
$ws.Range("D10").Value = 1.03050854360486
$ws.Range("E10").Value = 1.016410214593027
$ws.Range("F10").Value = 1.010667792754732
$ws.Range("I10").Value = 1.031001184960753
$ws.Range("J10").Value = 1.021040449966748
$ws.Range("K10").Value = 1.033973877930194
$ws.Range("L10").Value = 1.019927296992984
$ws.Range("M10").Value = 1.014206410942817
$ws.Range("N10").Value = 1.011103105670294
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.013984850312608
$ws.Range("D11").Value = 1.030184526244864
$ws.Range("E11").Value = 1.015928645792188
$ws.Range("F11").Value = 1.009667449515966
$ws.Range("I11").Value = 1.030906215921613
$ws.Range("J11").Value = 1.020707923933341
$ws.Range("K11").Value = 1.033773323468578
$ws.Range("L11").Value = 1.019571787669577
$ws.Range("M11").Value = 1.01333499829042
$ws.Range("N11").Value = 1.010991758143142
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.013772356623113
$ws.Range("D12").Value = 1.03006419486757
$ws.Range("E12").Value = 1.015749877935133
$ws.Range("F12").Value = 1.009295878735913
$ws.Range("I12").Value = 1.030870730752748
$ws.Range("J12").Value = 1.02058433300238
$ws.Range("K12").Value = 1.033698715922452
$ws.Range("L12").Value = 1.019439719548922
$ws.Range("M12").Value = 1.013011233792705
$ws.Range("N12").Value = 1.010950365000452
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.013817935510125
$ws.Range("D13").Value = 1.030090005284325
$ws.Range("E13").Value = 1.015788219336374
$ws.Range("F13").Value = 1.009375581957831
$ws.Range("I13").Value = 1.030878351912572
$ws.Range("J13").Value = 1.020610847110551
$ws.Range("K13").Value = 1.033714724577073
$ws.Range("L13").Value = 1.019468049332292
$ws.Range("M13").Value = 1.013080686205478
$ws.Range("N13").Value = 1.010959245493667
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.013967284738225
$ws.Range("D14").Value = 1.030174579134795
$ws.Range("E14").Value = 1.015913866572673
$ws.Range("F14").Value = 1.009636735364475
$ws.Range("I14").Value = 1.030903286976488
$ws.Range("J14").Value = 1.02069770941675
$ws.Range("K14").Value = 1.033767158678348
$ws.Range("L14").Value = 1.019560871204136
$ws.Range("M14").Value = 1.013308237535683
$ws.Range("N14").Value = 1.010988337260074
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.014059308789047
$ws.Range("D15").Value = 1.030226690986185
$ws.Range("E15").Value = 1.015991296348523
$ws.Range("F15").Value = 1.009797640604566
$ws.Range("I15").Value = 1.030918622552496
$ws.Range("J15").Value = 1.020751218091515
$ws.Range("K15").Value = 1.033799450157405
$ws.Range("L15").Value = 1.019618059698134
$ws.Range("M15").Value = 1.013448428284135
$ws.Range("N15").Value = 1.0110062571939
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.014595005577131
$ws.Range("D16").Value = 1.030530050716699
$ws.Range("E16").Value = 1.016442189862692
$ws.Range("F16").Value = 1.010734182702766
$ws.Range("I16").Value = 1.031007458348545
$ws.Range("J16").Value = 1.021062507961469
$ws.Range("K16").Value = 1.033987172209001
$ws.Range("L16").Value = 1.019950888691154
$ws.Range("M16").Value = 1.014264232329044
$ws.Range("N16").Value = 1.01111049070426
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.014931101600696
$ws.Range("D17").Value = 1.030720379610115
$ws.Range("E17").Value = 1.016725215461452
$ws.Range("F17").Value = 1.011321658673417
$ws.Range("I17").Value = 1.031062809007082
$ws.Range("J17").Value = 1.021257636297708
$ws.Range("K17").Value = 1.034104723355159
$ws.Range("L17").Value = 1.020159634333007
$ws.Range("M17").Value = 1.014775821449051
$ws.Range("N17").Value = 1.011175813418596
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.015127164971684
$ws.Range("D18").Value = 1.0308314088721
$ws.Range("E18").Value = 1.016890368440957
$ws.Range("F18").Value = 1.011664327982812
$ws.Range("I18").Value = 1.031094959298497
$ws.Range("J18").Value = 1.02137140242946
$ws.Range("K18").Value = 1.03417321583841
$ws.Range("L18").Value = 1.020281381478436
$ws.Range("M18").Value = 1.015074172541863
$ws.Range("N18").Value = 1.011213893312187
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.015194021616819
$ws.Range("D19").Value = 1.030869269243113
$ws.Range("E19").Value = 1.016946693007001
$ws.Range("F19").Value = 1.011781170287026
$ws.Range("I19").Value = 1.031105898837543
$ws.Range("J19").Value = 1.021410185435503
$ws.Range("K19").Value = 1.03419655757235
$ws.Range("L19").Value = 1.020322892288549
$ws.Range("M19").Value = 1.015175894176703
$ws.Range("N19").Value = 1.011226873880401
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.014895039167662
$ws.Range("D20").Value = 1.030699957698219
$ws.Range("E20").Value = 1.01669484235691
$ws.Range("F20").Value = 1.011258627631874
$ws.Range("I20").Value = 1.031056884346276
$ws.Range("J20").Value = 1.021236705932938
$ws.Range("K20").Value = 1.034092118784371
$ws.Range("L20").Value = 1.020137239000261
$ws.Range("M20").Value = 1.014720938000535
$ws.Range("N20").Value = 1.011168807157342
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.013923304057795
$ws.Range("D21").Value = 1.030149673592701
$ws.Range("E21").Value = 1.015876863598888
$ws.Range("F21").Value = 1.00955983220327
$ws.Range("I21").Value = 1.030895949999741
$ws.Range("J21").Value = 1.020672132727704
$ws.Range("K21").Value = 1.033751721244824
$ws.Range("L21").Value = 1.01953353791829
$ws.Range("M21").Value = 1.013241231668435
$ws.Range("N21").Value = 1.010979771389611
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.0133125563823
$ws.Range("D22").Value = 1.02980382259464
$ws.Range("E22").Value = 1.015363194513752
$ws.Range("F22").Value = 1.008491736654306
$ws.Range("I22").Value = 1.030793552650656
$ws.Range("J22").Value = 1.020316723783737
$ws.Range("K22").Value = 1.033537048328453
$ws.Range("L22").Value = 1.019153873974862
$ws.Range("M22").Value = 1.012310398926947
$ws.Range("N22").Value = 1.010860722105623
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.013636304243865
$ws.Range("D23").Value = 1.029987151447028
$ws.Range("E23").Value = 1.01563544054852
$ws.Range("F23").Value = 1.009057955636672
$ws.Range("I23").Value = 1.030847950118755
$ws.Range("J23").Value = 1.020505174368621
$ws.Range("K23").Value = 1.033650911907876
$ws.Range("L23").Value = 1.019355149701811
$ws.Range("M23").Value = 1.012803898207389
$ws.Range("N23").Value = 1.010923850825135
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.01491133415388
$ws.Range("D24").Value = 1.030709185439217
$ws.Range("E24").Value = 1.016708566442492
$ws.Range("F24").Value = 1.011287108631219
$ws.Range("I24").Value = 1.031059561862045
$ws.Range("J24").Value = 1.021246163614955
$ws.Range("K24").Value = 1.03409781447337
$ws.Range("L24").Value = 1.02014735852004
$ws.Range("M24").Value = 1.014745737622482
$ws.Range("N24").Value = 1.011171973052386
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.016391885220921
$ws.Range("D25").Value = 1.031547594543803
$ws.Range("E25").Value = 1.017956570505481
$ws.Range("F25").Value = 1.013874112181892
$ws.Range("I25").Value = 1.031299835538805
$ws.Range("J25").Value = 1.022104120309258
$ws.Range("K25").Value = 1.034613548173763
$ws.Range("L25").Value = 1.021066253147438
$ws.Range("M25").Value = 1.01699718109482
$ws.Range("N25").Value = 1.011459051760254
